$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style index 1: bold, centered, bordered) from an existing
# column-A data cell down into the new rows so the appended cells match the
# existing look-and-feel of column A.
$ws.Range("A205").Copy()
$ws.Range("A206:A217").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Append the new rows (206-217) with sequence numbers 204-215 in column A
# and their corresponding normalized values in column B.
$ws.Cells.Item(206, 1).Value = 204
$ws.Cells.Item(206, 2).Value = [double]"3.700743415417188E-17"

$ws.Cells.Item(207, 1).Value = 205
$ws.Cells.Item(207, 2).Value = 0

$ws.Cells.Item(208, 1).Value = 206
$ws.Cells.Item(208, 2).Value = [double]"6.245004513516506E-18"

$ws.Cells.Item(209, 1).Value = 207
$ws.Cells.Item(209, 2).Value = [double]"4.317533984653387E-17"

$ws.Cells.Item(210, 1).Value = 208
$ws.Cells.Item(210, 2).Value = [double]"-6.938893903907228E-18"

$ws.Cells.Item(211, 1).Value = 209
$ws.Cells.Item(211, 2).Value = [double]"6.344131569286608E-17"

$ws.Cells.Item(212, 1).Value = 210
$ws.Cells.Item(212, 2).Value = [double]"4.394632805807911E-17"

$ws.Cells.Item(213, 1).Value = 211
$ws.Cells.Item(213, 2).Value = [double]"-2.775557561562892E-18"

$ws.Cells.Item(214, 1).Value = 212
$ws.Cells.Item(214, 2).Value = [double]"-3.469446951953614E-18"

$ws.Cells.Item(215, 1).Value = 213
$ws.Cells.Item(215, 2).Value = [double]"4.625929269271486E-17"

$ws.Cells.Item(216, 1).Value = 214
$ws.Cells.Item(216, 2).Value = 0

$ws.Cells.Item(217, 1).Value = 215
$ws.Cells.Item(217, 2).Value = 0
